$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.320746
$ws.Range("H2").Value = 288.962238
$ws.Range("I2").Value = 0.3809824610908788
$ws.Range("J2").Value = 0.3809824610908788
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 20.424575
$ws.Range("N2").Value = 61.273725
$ws.Range("O2").Value = 0.203732656096709
$ws.Range("P2").Value = 0.2037326560967089
$ws.Range("Q2").Value = 1967.31030073295
$ws.Range("R2").Value = 17705.79270659655
$ws.Range("S2").Value = 0.07761856872430581
$ws.Range("T2").Value = 0.07761856872430579

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.320746
$ws.Range("H3").Value = 288.962238
$ws.Range("I3").Value = 0.3809824610908788
$ws.Range("J3").Value = 0.3809824610908788
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 57.16769933333334
$ws.Range("N3").Value = 171.503098
$ws.Range("O3").Value = 0.5702408607336045
$ws.Range("P3").Value = 0.5702408607336045
$ws.Range("Q3").Value = 5506.43544689037
$ws.Range("R3").Value = 49557.91902201333
$ws.Range("S3").Value = 0.2172517665368697
$ws.Range("T3").Value = 0.2172517665368697

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.320746
$ws.Range("H4").Value = 288.962238
$ws.Range("I4").Value = 0.3809824610908788
$ws.Range("J4").Value = 0.3809824610908788
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1182046666666667
$ws.Range("N4").Value = 0.354614
$ws.Range("O4").Value = 0.001179077200040937
$ws.Range("P4").Value = 0.001179077200040937
$ws.Range("Q4").Value = 11.38556167401467
$ws.Range("R4").Value = 102.470055066132
$ws.Range("S4").Value = 0.0004492077334877385
$ws.Range("T4").Value = 0.0004492077334877385

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 96.320746
$ws.Range("H5").Value = 288.962238
$ws.Range("I5").Value = 0.3809824610908788
$ws.Range("J5").Value = 0.3809824610908788
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 22.54136766666667
$ws.Range("N5").Value = 67.62410300000001
$ws.Range("O5").Value = 0.2248474059696456
$ws.Range("P5").Value = 0.2248474059696456
$ws.Range("Q5").Value = 2171.201349513613
$ws.Range("R5").Value = 19540.81214562252
$ws.Range("S5").Value = 0.08566291809621555
$ws.Range("T5").Value = 0.08566291809621554

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.46467533333333
$ws.Range("H6").Value = 55.394026
$ws.Range("I6").Value = 0.07303429161291354
$ws.Range("J6").Value = 0.07303429161291354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 20.424575
$ws.Range("N6").Value = 61.273725
$ws.Range("O6").Value = 0.203732656096709
$ws.Range("P6").Value = 0.2037326560967089
$ws.Range("Q6").Value = 377.1331461963167
$ws.Range("R6").Value = 3394.19831576685
$ws.Range("S6").Value = 0.01487947021644047
$ws.Range("T6").Value = 0.01487947021644047

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.46467533333333
$ws.Range("H7").Value = 55.394026
$ws.Range("I7").Value = 0.07303429161291354
$ws.Range("J7").Value = 0.07303429161291354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 57.16769933333334
$ws.Range("N7").Value = 171.503098
$ws.Range("O7").Value = 0.5702408607336045
$ws.Range("P7").Value = 0.5702408607336045
$ws.Range("Q7").Value = 1055.583007743616
$ws.Range("R7").Value = 9500.247069692548
$ws.Range("S7").Value = 0.04164713731241688
$ws.Range("T7").Value = 0.04164713731241688

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.46467533333333
$ws.Range("H8").Value = 55.394026
$ws.Range("I8").Value = 0.07303429161291354
$ws.Range("J8").Value = 0.07303429161291354
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.1182046666666667
$ws.Range("N8").Value = 0.354614
$ws.Range("O8").Value = 0.001179077200040937
$ws.Range("P8").Value = 0.001179077200040937
$ws.Range("Q8").Value = 2.182610792884889
$ws.Range("R8").Value = 19.643497135964
$ws.Range("S8").Value = 0.00008611306806192736
$ws.Range("T8").Value = 0.00008611306806192736

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.46467533333333
$ws.Range("H9").Value = 55.394026
$ws.Range("I9").Value = 0.07303429161291354
$ws.Range("J9").Value = 0.07303429161291354
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 22.54136766666667
$ws.Range("N9").Value = 67.62410300000001
$ws.Range("O9").Value = 0.2248474059696456
$ws.Range("P9").Value = 0.2248474059696456
$ws.Range("Q9").Value = 416.2190355342976
$ws.Range("R9").Value = 3745.971319808678
$ws.Range("S9").Value = 0.01642157101599426
$ws.Range("T9").Value = 0.01642157101599425

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 122.909391
$ws.Range("H10").Value = 368.728173
$ws.Range("I10").Value = 0.4861499128584522
$ws.Range("J10").Value = 0.4861499128584522
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.424575
$ws.Range("N10").Value = 61.273725
$ws.Range("O10").Value = 0.203732656096709
$ws.Range("P10").Value = 0.2037326560967089
$ws.Range("Q10").Value = 2510.372074683825
$ws.Range("R10").Value = 22593.34867215442
$ws.Range("S10").Value = 0.09904461300783607
$ws.Range("T10").Value = 0.09904461300783605

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 122.909391
$ws.Range("H11").Value = 368.728173
$ws.Range("I11").Value = 0.4861499128584522
$ws.Range("J11").Value = 0.4861499128584522
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 57.16769933333334
$ws.Range("N11").Value = 171.503098
$ws.Range("O11").Value = 0.5702408607336045
$ws.Range("P11").Value = 0.5702408607336045
$ws.Range("Q11").Value = 7026.447109931106
$ws.Range("R11").Value = 63238.02398937996
$ws.Range("S11").Value = 0.2772225447539706
$ws.Range("T11").Value = 0.2772225447539706

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 122.909391
$ws.Range("H12").Value = 368.728173
$ws.Range("I12").Value = 0.4861499128584522
$ws.Range("J12").Value = 0.4861499128584522
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.1182046666666667
$ws.Range("N12").Value = 0.354614
$ws.Range("O12").Value = 0.001179077200040937
$ws.Range("P12").Value = 0.001179077200040937
$ws.Range("Q12").Value = 14.528463593358
$ws.Range("R12").Value = 130.756172340222
$ws.Range("S12").Value = 0.0005732082780532891
$ws.Range("T12").Value = 0.0005732082780532891

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 122.909391
$ws.Range("H13").Value = 368.728173
$ws.Range("I13").Value = 0.4861499128584522
$ws.Range("J13").Value = 0.4861499128584522
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 22.54136766666667
$ws.Range("N13").Value = 67.62410300000001
$ws.Range("O13").Value = 0.2248474059696456
$ws.Range("P13").Value = 0.2248474059696456
$ws.Range("Q13").Value = 2770.545772217091
$ws.Range("R13").Value = 24934.91194995382
$ws.Range("S13").Value = 0.1093095468185922
$ws.Range("T13").Value = 0.1093095468185922

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.127183
$ws.Range("H14").Value = 45.381549
$ws.Range("I14").Value = 0.05983333443775553
$ws.Range("J14").Value = 0.05983333443775553
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 20.424575
$ws.Range("N14").Value = 61.273725
$ws.Range("O14").Value = 0.203732656096709
$ws.Range("P14").Value = 0.2037326560967089
$ws.Range("Q14").Value = 308.966283722225
$ws.Range("R14").Value = 2780.696553500025
$ws.Range("S14").Value = 0.01219000414812662
$ws.Range("T14").Value = 0.01219000414812662

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.127183
$ws.Range("H15").Value = 45.381549
$ws.Range("I15").Value = 0.05983333443775553
$ws.Range("J15").Value = 0.05983333443775553
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 57.16769933333334
$ws.Range("N15").Value = 171.503098
$ws.Range("O15").Value = 0.5702408607336045
$ws.Range("P15").Value = 0.5702408607336045
$ws.Range("Q15").Value = 864.7862495043114
$ws.Range("R15").Value = 7783.076245538803
$ws.Range("S15").Value = 0.03411941213034733
$ws.Range("T15").Value = 0.03411941213034733

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.127183
$ws.Range("H16").Value = 45.381549
$ws.Range("I16").Value = 0.05983333443775553
$ws.Range("J16").Value = 0.05983333443775553
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.1182046666666667
$ws.Range("N16").Value = 0.354614
$ws.Range("O16").Value = 0.001179077200040937
$ws.Range("P16").Value = 0.001179077200040937
$ws.Range("Q16").Value = 1.788103624120667
$ws.Range("R16").Value = 16.092932617086
$ws.Range("S16").Value = 0.00007054812043798173
$ws.Range("T16").Value = 0.00007054812043798173

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.127183
$ws.Range("H17").Value = 45.381549
$ws.Range("I17").Value = 0.05983333443775553
$ws.Range("J17").Value = 0.05983333443775553
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 22.54136766666667
$ws.Range("N17").Value = 67.62410300000001
$ws.Range("O17").Value = 0.2248474059696456
$ws.Range("P17").Value = 0.2248474059696456
$ws.Range("Q17").Value = 340.9873937639497
$ws.Range("R17").Value = 340.9873937639497
$ws.Range("S17").Value = 0.0134533700388436
$ws.Range("T17").Value = 0.0134533700388436
